# "Generate Report for Handoff" - append a new handoff row for
# a0f9f52f-6c78-467f-b9c4-0297c6599336 to all three sheets
# (Overview, zh-cn, de-de), mirroring the existing rows' shape,
# hyperlinks and formatting.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$fileBase  = "a0f9f52f-6c78-467f-b9c4-0297c6599336"
$mdName    = "$fileBase.md"
$zhName    = "$fileBase.9f35d8939bd5bffcc023ebcad41537e9347c5859.zh-cn.xlf"
$deName    = "$fileBase.9f35d8939bd5bffcc023ebcad41537e9347c5859.de-de.xlf"

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/65322263e2318b89ddbb50ebe95ebd90954c6123/e2e/$mdName"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/23caed3dd6c55d224000671cb08f43707e11c633/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$zhName"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c91d31de419d62bbd3008c787323a642be0f14d/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$deName"

$statusText  = "Ready for handoff"
$includeText = "Include"
$mdExt       = ".md"
$epoch       = "0001-01-01 00:00:00"

$mdHandoffDatetime = "2016-03-20 16:50:47"
$zhHandoffDatetime = "2016-03-20 16:50:39"

# ---------------------------------------------------------------
# Sheet 1: Overview  (row 7 -> A:D)
# ---------------------------------------------------------------
$ws1.Hyperlinks.Add($ws1.Cells.Item(7, 1), $mdUrl, "", "", $mdName) | Out-Null
$ws1.Cells.Item(7, 1).Style = "Normal"

$ws1.Cells.Item(7, 2).Value = $statusText
$ws1.Cells.Item(7, 2).Style = "Normal"

$ws1.Cells.Item(7, 3).Value = $statusText
$ws1.Cells.Item(7, 3).Style = "Normal"

$ws1.Cells.Item(7, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws1.Cells.Item(7, 4).Value = $mdHandoffDatetime

# ---------------------------------------------------------------
# Sheet 2: zh-cn  (row 7 -> A,B,C,D,E,H,J)
# ---------------------------------------------------------------
$ws2.Hyperlinks.Add($ws2.Cells.Item(7, 1), $mdUrl, "", "", $mdName) | Out-Null
$ws2.Cells.Item(7, 1).Style = "Normal"

$ws2.Cells.Item(7, 2).Value = $mdExt
$ws2.Cells.Item(7, 2).Style = "Normal"

$ws2.Cells.Item(7, 3).Value = $statusText
$ws2.Cells.Item(7, 3).Style = "Normal"

$ws2.Hyperlinks.Add($ws2.Cells.Item(7, 4), $zhUrl, "", "", $zhName) | Out-Null
$ws2.Cells.Item(7, 4).Style = "Normal"

$ws2.Cells.Item(7, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(7, 5).Value = $zhHandoffDatetime

$ws2.Cells.Item(7, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(7, 8).Value = $epoch

$ws2.Cells.Item(7, 10).Value = $includeText
$ws2.Cells.Item(7, 10).Style = "Normal"

# ---------------------------------------------------------------
# Sheet 3: de-de  (row 7 -> A,B,C,D,E,H,J)
# ---------------------------------------------------------------
$ws3.Hyperlinks.Add($ws3.Cells.Item(7, 1), $mdUrl, "", "", $mdName) | Out-Null
$ws3.Cells.Item(7, 1).Style = "Normal"

$ws3.Cells.Item(7, 2).Value = $mdExt
$ws3.Cells.Item(7, 2).Style = "Normal"

$ws3.Cells.Item(7, 3).Value = $statusText
$ws3.Cells.Item(7, 3).Style = "Normal"

$ws3.Hyperlinks.Add($ws3.Cells.Item(7, 4), $deUrl, "", "", $deName) | Out-Null
$ws3.Cells.Item(7, 4).Style = "Normal"

$ws3.Cells.Item(7, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(7, 5).Value = $mdHandoffDatetime

$ws3.Cells.Item(7, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(7, 8).Value = $epoch

$ws3.Cells.Item(7, 10).Value = $includeText
$ws3.Cells.Item(7, 10).Style = "Normal"
